$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.080.41'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '1.871.64'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '312.84'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = '0.5134'
$ws.Range("E7").Value = '  +1.77%  '
$ws.Range("D8").Value = '0.3877'
$ws.Range("E8").Value = '  +1.25%  '
$ws.Range("D9").Value = '0.08365'
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("D11").Value = '41.46'
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").Value = '6.193'
$ws.Range("E12").Value = '  -1.26%  '
$ws.Range("D13").Value = '20.55'
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").Value = '1.863.79'
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("D15").Value = '7.289'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '0.00001104'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '90.98'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").Value = '0.06657'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("E20").Value = '  -2.40%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("E22").Value = '  -1.47%  '
$ws.Range("D23").Value = '28.114.81'
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("D24").Value = '11.11'
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("D25").Value = '2.249'
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("D26").Value = '2.080.81'
$ws.Range("E26").Value = '  -0.84%  '
$ws.Range("D27").Value = '2.482'
$ws.Range("E27").Value = '  -4.36%  '
$ws.Range("D28").Value = '158.28'
$ws.Range("E28").Value = '  +1.22%  '
$ws.Range("D29").Value = '20.57'
$ws.Range("E29").Value = '  -0.76%  '
$ws.Range("D30").Value = '124.93'
$ws.Range("E30").Value = '  -1.54%  '
$ws.Range("E31").Value = '  +0.96%  '
$ws.Range("D32").Value = '1.036'
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").Value = '5.885'
$ws.Range("E33").Value = '  +4.30%  '
$ws.Range("D34").Value = '3.593'
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("E35").Value = '  -3.19%  '
$ws.Range("D36").Value = '0.02434'
$ws.Range("E36").Value = '  -0.81%  '
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("D39").Value = '1.204'
$ws.Range("E39").Value = '  -2.52%  '
$ws.Range("D40").Value = '0.6480'
$ws.Range("E40").Value = '  +1.48%  '
$ws.Range("D41").Value = '5.009'
$ws.Range("E41").Value = '  +2.27%  '
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("E43").Value = '  -0.75%  '
$ws.Range("D44").Value = '0.6093'
$ws.Range("E44").Value = '  +0.48%  '
$ws.Range("D45").Value = '13.02'
$ws.Range("E45").Value = '  -1.20%  '
$ws.Range("D46").Value = '3.679'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("D48").Value = '2.011'
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").Value = '121.12'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = '78.01'
$ws.Range("E51").Value = '  -3.29%  '
